# Attendance sheet cleanup:
#  - remove the trailing placeholder/duplicate rows (rows 7-11)
#  - fix up roll numbers in the remaining rows, including marking
#    previously-unresolved rows with the correct roll number + status

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unneeded rows 7:11 (duplicates / unresolved N/A placeholders)
$ws.Rows("7:11").Delete()

# Row 2: was an unresolved "N/A"/"Unknown" row with no status -> roll 21115030, Present
$ws.Range("A2").Value = "'21115030"
$ws.Range("B2").Value = "'21115030"
$ws.Range("D2").Value = "P"

# Row 3: roll number correction 21115024 -> 21115021
$ws.Range("A3").Value = "'21115021"
$ws.Range("B3").Value = "'21115021"

# Row 4: roll number correction 21115086 -> 21115024
$ws.Range("A4").Value = "'21115024"
$ws.Range("B4").Value = "'21115024"

# Row 5: roll number correction 21115024 -> 21115086
$ws.Range("A5").Value = "'21115086"
$ws.Range("B5").Value = "'21115086"

# Row 6: roll number correction 21115038 -> 21115080
$ws.Range("A6").Value = "'21115080"
$ws.Range("B6").Value = "'21115080"
